$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the price for row 13 (resistor networks replacing discrete parts)
$ws.Range("C13").Value = 0.0156

# Add part number / LCSC PN references for the new 1206 resistor networks
$ws.Range("E13").Value = "BZT52C3V6"
$ws.Range("F13").Value = "C173412"

# Match the cell that was left selected after the edit
$ws.Range("F13").Select()
